$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 9920216
$ws.Range("I8").Value = 11022240
$ws.Range("K8").Value = 33066720
$ws.Range("M8").Value = -33066581

$ws.Range("H31").Value = 542.8570999999999
$ws.Range("I31").Value = 466.66666
$ws.Range("K31").Value = 1399.99998
$ws.Range("M31").Value = -1169.99998

$ws.Range("H40").Value = 2557.1428
$ws.Range("J40").Value = 7000
$ws.Range("L40").Value = 7000
$ws.Range("N40").Value = -7350

$ws.Range("H64").Value = 3862.926
$ws.Range("I64").Value = 3531.96
$ws.Range("K64").Value = 3531.96
$ws.Range("M64").Value = -3283.96

$ws.Range("H67").Value = 3862.926
$ws.Range("I67").Value = 3531.96
$ws.Range("K67").Value = 3531.96
$ws.Range("M67").Value = -2673.96

$ws.Range("H74").Value = 3893.4707
$ws.Range("I74").Value = 2773.625
$ws.Range("J74").Value = 4888.8887
$ws.Range("K74").Value = 2773.625
$ws.Range("L74").Value = 4888.8887
$ws.Range("M74").Value = -1837.625
$ws.Range("N74").Value = -6760.8887

$ws.Range("H77").Value = 3893.4707
$ws.Range("I77").Value = 2773.625
$ws.Range("J77").Value = 4888.8887
$ws.Range("K77").Value = 13868.125
$ws.Range("L77").Value = 24444.4435
$ws.Range("M77").Value = -9188.125
$ws.Range("N77").Value = -33804.4435

$ws.Range("H137").Value = 1053.1072
$ws.Range("I137").Value = 990.45
$ws.Range("J137").Value = 1209.75
$ws.Range("K137").Value = 2971.35
$ws.Range("L137").Value = 3629.25
$ws.Range("M137").Value = -421.3500000000004
$ws.Range("N137").Value = -8729.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1857.6666
$ws.Range("I74").Value = 1779.2
$ws.Range("K74").Value = 1779.2
$ws.Range("M74").Value = -905.2

$ws.Range("H77").Value = 1857.6666
$ws.Range("I77").Value = 1779.2
$ws.Range("K77").Value = 8896
$ws.Range("M77").Value = -4528

$ws.Range("H82").Value = 12672.5
$ws.Range("I82").Value = 164
$ws.Range("K82").Value = 164
$ws.Range("M82").Value = 197

$ws.Range("H85").Value = 12672.5
$ws.Range("I85").Value = 164
$ws.Range("K85").Value = 164
$ws.Range("M85").Value = 1084

$ws.Range("H132").Value = 1728.3158
$ws.Range("I132").Value = 1660.5
$ws.Range("K132").Value = 4981.5
$ws.Range("M132").Value = -2451.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 713.25
$ws.Range("I7").Value = 651.5
$ws.Range("J7").Value = 775
$ws.Range("K7").Value = 651.5
$ws.Range("L7").Value = 775
$ws.Range("M7").Value = -538.5
$ws.Range("N7").Value = -1001

$ws.Range("H86").Value = 2341.2632
$ws.Range("I86").Value = 2436.6365
$ws.Range("K86").Value = 2436.6365
$ws.Range("M86").Value = -1313.6365

$ws.Range("H89").Value = 2341.2632
$ws.Range("I89").Value = 2436.6365
$ws.Range("K89").Value = 12183.1825
$ws.Range("M89").Value = -6567.182500000001

$ws.Range("H132").Value = 130913.2
$ws.Range("J132").Value = 130913.2
$ws.Range("L132").Value = 130913.2
$ws.Range("N132").Value = -141033.2

$ws.Range("H134").Value = 1958.9333
$ws.Range("I134").Value = 1721.1111
$ws.Range("K134").Value = 5163.3333
$ws.Range("M134").Value = -2628.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12857.473
$ws.Range("I31").Value = 3501.9524
$ws.Range("K31").Value = 3501.9524
$ws.Range("M31").Value = -3206.9524

$ws.Range("H34").Value = 12857.473
$ws.Range("I34").Value = 3501.9524
$ws.Range("K34").Value = 3501.9524
$ws.Range("M34").Value = -3299.9524

$ws.Range("H35").Value = 1899.7142
$ws.Range("I35").Value = 566.6667
$ws.Range("K35").Value = 566.6667
$ws.Range("M35").Value = -272.6667

$ws.Range("H109").Value = 84000
$ws.Range("J109").Value = 84000
$ws.Range("L109").Value = 84000
$ws.Range("N109").Value = -86080

$ws.Range("H132").Value = 2952.4443
$ws.Range("I132").Value = 2856.5938
$ws.Range("K132").Value = 8569.7814
$ws.Range("M132").Value = -6039.7814

$ws.Range("H134").Value = 2834.5
$ws.Range("J134").Value = 5822
$ws.Range("L134").Value = 17466
$ws.Range("N134").Value = -22536

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4916381.5
$ws.Range("J4").Value = 4882533.5
$ws.Range("L4").Value = 14647600.5
$ws.Range("N4").Value = -14647824.5

$ws.Range("H14").Value = 1309.1875
$ws.Range("I14").Value = 1309.1875
$ws.Range("K14").Value = 3927.5625
$ws.Range("M14").Value = -3754.5625

$ws.Range("H21").Value = 783.1667
$ws.Range("I21").Value = 759.6
$ws.Range("J21").Value = 901
$ws.Range("K21").Value = 2278.8
$ws.Range("L21").Value = 2703
$ws.Range("M21").Value = -2105.8
$ws.Range("N21").Value = -3049

$ws.Range("H24").Value = 6700
$ws.Range("J24").Value = 10000
$ws.Range("L24").Value = 30000
$ws.Range("N24").Value = -30460

$ws.Range("H25").Value = 262.5
$ws.Range("I25").Value = 225
$ws.Range("K25").Value = 675
$ws.Range("M25").Value = -506

$ws.Range("H30").Value = 262.5
$ws.Range("I30").Value = 225
$ws.Range("K30").Value = 675
$ws.Range("M30").Value = -573

$ws.Range("H43").Value = 5000
$ws.Range("J43").Value = 5000
$ws.Range("L43").Value = 15000
$ws.Range("N43").Value = -15228

$ws.Range("H48").Value = 1200
$ws.Range("J48").Value = 0
$ws.Range("L48").Value = 0
$ws.Range("N48").ClearContents()

$ws.Range("H54").Value = 8333
$ws.Range("J54").Value = 8333
$ws.Range("L54").Value = 24999
$ws.Range("N54").Value = -26117

$ws.Range("H64").Value = 1900
$ws.Range("I64").Value = 0
$ws.Range("K64").Value = 0
$ws.Range("M64").ClearContents()

$ws.Range("H67").Value = 1900
$ws.Range("I67").Value = 0
$ws.Range("K67").Value = 0
$ws.Range("M67").ClearContents()

$ws.Range("H88").Value = 10552.4
$ws.Range("J88").Value = 10552.4
$ws.Range("L88").Value = 31657.2
$ws.Range("N88").Value = -32513.2

$ws.Range("H91").Value = 10552.4
$ws.Range("J91").Value = 10552.4
$ws.Range("L91").Value = 31657.2
$ws.Range("N91").Value = -34621.2

$ws.Range("H94").Value = 6380
$ws.Range("J94").Value = 7375
$ws.Range("L94").Value = 22125
$ws.Range("N94").Value = -23477

$ws.Range("H105").Value = 14998
$ws.Range("J105").Value = 14998
$ws.Range("L105").Value = 44994
$ws.Range("N105").Value = -50236

$ws.Range("H106").Value = 5825.4287
$ws.Range("J106").Value = 5916.75
$ws.Range("L106").Value = 17750.25
$ws.Range("N106").Value = -19642.25

$ws.Range("H120").Value = 5423.8887
$ws.Range("I120").Value = 5423.8887
$ws.Range("J120").Value = 0
$ws.Range("K120").Value = 16271.6661
$ws.Range("L120").Value = 0
$ws.Range("M120").Value = -11433.6661
$ws.Range("N120").ClearContents()

$ws.Range("H126").Value = 12500
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 12500
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 37500
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -47380

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()

$ws.Range("H63").Value = 23114
$ws.Range("J63").Value = 23114
$ws.Range("L63").Value = 23114
$ws.Range("N63").Value = -24486

$ws.Range("H66").Value = 23114
$ws.Range("J66").Value = 23114
$ws.Range("L66").Value = 69342
$ws.Range("N66").Value = -76206

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1603.25
$ws.Range("I22").Value = 1026.6666
$ws.Range("J22").Value = 3333
$ws.Range("K22").Value = 1026.6666
$ws.Range("L22").Value = 3333
$ws.Range("M22").Value = -731.6666
$ws.Range("N22").Value = -3923

$ws.Range("H27").Value = 1603.25
$ws.Range("I27").Value = 1026.6666
$ws.Range("J27").Value = 3333
$ws.Range("K27").Value = 1026.6666
$ws.Range("L27").Value = 3333
$ws.Range("M27").Value = -919.6666
$ws.Range("N27").Value = -3547

$ws.Range("H46").Value = 25841.777
$ws.Range("I46").Value = 54275.125
$ws.Range("J46").Value = 3095.1
$ws.Range("K46").Value = 54275.125
$ws.Range("L46").Value = 3095.1
$ws.Range("M46").Value = -54087.125
$ws.Range("N46").Value = -3471.1

$ws.Range("H55").Value = 397.13635
$ws.Range("I55").Value = 452.15384
$ws.Range("J55").Value = 317.66666
$ws.Range("K55").Value = 452.15384
$ws.Range("L55").Value = 317.66666
$ws.Range("M55").Value = -279.15384
$ws.Range("N55").Value = -663.66666

$ws.Range("H93").Value = 10340.595
$ws.Range("I93").Value = 1368.5834
$ws.Range("J93").Value = 333333
$ws.Range("K93").Value = 1368.5834
$ws.Range("L93").Value = 333333
$ws.Range("M93").Value = -120.5834
$ws.Range("N93").Value = -335829

$ws.Range("H132").Value = 4160.278
$ws.Range("I132").Value = 3926
$ws.Range("K132").Value = 11778
$ws.Range("M132").Value = -9248

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1724.2222
$ws.Range("I132").Value = 1113
$ws.Range("K132").Value = 3339
$ws.Range("M132").Value = -809
